$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Binary Search 1")
$ws3 = $wb.Worksheets.Item("Notes")

# New row 11 values, entered in an order that reproduces the expected
# shared-string insertion order (Question, then Page No., then link text).
$ws1.Range("D11").Value = "Single elements in a sorted array"
$ws1.Range("C11").Value = "Binary S1 26"
$ws1.Range("B11").Value = 9

# Apply the wrap/left/top style used for the new "Question" cell (D11).
$ws1.Range("D11").HorizontalAlignment = -4131   # xlLeft
$ws1.Range("D11").VerticalAlignment = -4160     # xlTop
$ws1.Range("D11").WrapText = $true

# The wrapped question text makes the row taller (two lines instead of one).
$ws1.Rows.Item(11).RowHeight = 28.8

# Add the hyperlink for E11; pass the URL as the display text so the
# "display" attribute mirrors the address like the other links in this
# sheet, then overwrite the cell text with the friendly title.
$ws1.Hyperlinks.Add($ws1.Range("E11"), "https://leetcode.com/problems/single-element-in-a-sorted-array/description/", [ref]$null, [ref]$null, "https://leetcode.com/problems/single-element-in-a-sorted-array/description/") | Out-Null
$ws1.Range("E11").Value = "Single Element in a Sorted Array - LeetCode"

# Match the existing "Hyperlink" cell style (as used on sheet "Notes")
# instead of the ad-hoc style Excel would otherwise synthesize.
$ws3.Range("F3").Copy()
$ws1.Range("E11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update selection to reflect where the edit left the cursor.
$ws1.Range("F11").Select()
